$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values in column A
$ws.Range("A1").Value = 0.07655054299873143
$ws.Range("A2").Value = -0.0099999998302813253
$ws.Range("A3").Value = -0.008999999832026262
$ws.Range("A4").Value = 0.28399655374833088
$ws.Range("A5").Value = -0.0059999998385684705
$ws.Range("A6").Value = -0.0059999998332749271
$ws.Range("A7").Value = -0.019999999799068746
$ws.Range("A8").Value = -0.0075494879613238197
$ws.Range("A9").Value = -0.0059999998313680081
$ws.Range("A10").Value = -0.0059999998307844749
$ws.Range("A11").Value = -0.0044999998344543712
$ws.Range("A12").Value = 0.00095913507567235712
$ws.Range("A13").Value = -0.0059999998304869351
$ws.Range("A14").Value = -0.011999999816010742
$ws.Range("A15").Value = -0.0059999998308395419
$ws.Range("A16").Value = -0.0059999998314292924
$ws.Range("A17").Value = 0.035315595160493096
$ws.Range("A18").Value = -0.0089999998253640356
$ws.Range("A19").Value = -0.0089999998326528718
$ws.Range("A20").Value = -0.0089999998311274254
$ws.Range("A21").Value = -0.0089999998308805118
$ws.Range("A22").Value = -0.0089999998307099816
$ws.Range("A23").Value = -0.0608334481480588
$ws.Range("A24").Value = -0.041999999746789918
$ws.Range("A25").Value = -0.041999999745350181
$ws.Range("A26").Value = -0.005999999832582148
$ws.Range("A27").Value = -0.0059999998316726533
$ws.Range("A28").Value = -0.0059999998280071409
$ws.Range("A29").Value = -0.011999999810777595
$ws.Range("A30").Value = -0.019999999790099032
$ws.Range("A31").Value = -0.014999999800336283
$ws.Range("A32").Value = -0.02099999978549949
$ws.Range("A33").Value = -0.0059999998216886397

# Update column A width (source stored width 15.42578125 -> 15.7109375;
# this runtime quantizes ColumnWidth to 1/6-character steps, so 15.666666666666666
# is the closest representable stored width to the target)
$ws.Columns.Item(1).ColumnWidth = 14.86
